$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44449, 0, 0, 0),
    @(44450, 0, 0, 0),
    @(44451, 0, 0, 0),
    @(44452, 1, 1, 83.40283569641367),
    @(44453, 0, 1, 83.40283569641367),
    @(44454, 0, 1, 83.40283569641367),
    @(44455, 0, 1, 83.40283569641367),
    @(44456, 3, 4, 333.6113427856547),
    @(44457, 0, 4, 333.6113427856547),
    @(44458, 0, 4, 333.6113427856547),
    @(44459, 0, 3, 250.208507089241)
)

$lastRow = 374
$startRow = $lastRow + 1
$endRow = $startRow + $data.Length - 1

# Copy formatting from the last existing row down into the new rows
$srcRange = $ws.Range("A$lastRow`:D$lastRow")
$dstRange = $ws.Range("A$startRow`:D$endRow")
$srcRange.Copy($dstRange)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
